$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 5; $row++) {
    $ws.Cells.Item($row, 7).Value = 0.4760219657335256
    $ws.Cells.Item($row, 8).Value = 0.998
}
